$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The script re-ran and produced a new day's result (2025-04-08). The "NA"
# value that used to sit in C42 (for 2025-04-07) is cleared out - C42 goes
# blank, like all the other "Rien ne nous concerne..." rows - and the "NA"
# reading is now recorded against the newly appended row instead.
$ws.Cells.Item(42, 3).Value = ""

# Append the new row for 2025-04-08.
$dateCell = $ws.Cells.Item(43, 1)
# Force text formatting on the date column first so Excel doesn't
# reinterpret the literal "2025-04-08" string as a date serial number,
# then restore the default "Normal" style so no stray formatting is left
# behind on the cell.
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-04-08"
$dateCell.Style = "Normal"

$ws.Cells.Item(43, 2).Value = "Rien ne nous concerne aujourd'hui !"
$ws.Cells.Item(43, 3).Value = "NA"
$ws.Cells.Item(43, 4).Value = 1
